# InitProperty.xlsx: "unify the conception of DataNode, DataTable, Entity."
#
# The sheet that used to model a generic "Property" table is renamed to
# "DataNode" so its name matches the unified DataNode/DataTable/Entity
# vocabulary used elsewhere in the data-config pipeline. No cell values,
# formulas, or shared strings change - this is purely a rename + the
# incidental view-state (selected cell) that Excel records when a user
# makes the edit interactively.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet: "Property1" -> "DataNode"
$ws.Name = "DataNode"

# Move the active selection to where the editor left off (matches the
# <selection activeCell="E50" .../> recorded in the saved view state).
$ws.Range("E50").Select()
